$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 641.75
$ws.Range("I19").Value = 564.5
$ws.Range("J19").Value = 667.5
$ws.Range("K19").Value = 564.5
$ws.Range("L19").Value = 667.5
$ws.Range("M19").Value = -389.5
$ws.Range("N19").Value = -1017.5

$ws.Range("H33").Value = 128.90909
$ws.Range("I33").Value = 70.588234
$ws.Range("K33").Value = 70.588234
$ws.Range("M33").Value = 158.411766

$ws.Range("H98").Value = 339954.62
$ws.Range("I98").Value = 430492.44
$ws.Range("J98").Value = 3671.4285
$ws.Range("K98").Value = 430492.44
$ws.Range("L98").Value = 3671.4285
$ws.Range("M98").Value = -428994.44
$ws.Range("N98").Value = -6667.4285

$ws.Range("H122").Value = 339954.62
$ws.Range("I122").Value = 430492.44
$ws.Range("J122").Value = 3671.4285
$ws.Range("K122").Value = 1291477.32
$ws.Range("L122").Value = 11014.2855
$ws.Range("M122").Value = -1289027.32
$ws.Range("N122").Value = -15914.2855

$ws.Range("H125").Value = 10192582
$ws.Range("I125").Value = 399.5
$ws.Range("J125").Value = 12457512
$ws.Range("K125").Value = 3595.5
$ws.Range("L125").Value = 112117608
$ws.Range("M125").Value = -1135.5
$ws.Range("N125").Value = -112122528

$ws.Range("H137").Value = 1188.3208
$ws.Range("I137").Value = 738.7143
$ws.Range("J137").Value = 1483.375
$ws.Range("K137").Value = 2216.1429
$ws.Range("L137").Value = 4450.125
$ws.Range("M137").Value = 333.8571000000002
$ws.Range("N137").Value = -9550.125

$ws.Range("H138").Value = 10419083
$ws.Range("I138").Value = 2294.6924
$ws.Range("J138").Value = 22729832
$ws.Range("K138").Value = 6884.0772
$ws.Range("L138").Value = 68189496
$ws.Range("M138").Value = -1744.0772
$ws.Range("N138").Value = -68199776

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1334.6415
$ws.Range("I132").Value = 974.93475
$ws.Range("J132").Value = 3698.4285
$ws.Range("K132").Value = 2924.80425
$ws.Range("L132").Value = 11095.2855
$ws.Range("M132").Value = -394.8042500000001
$ws.Range("N132").Value = -16155.2855

$ws.Range("H133").Value = 44563.355
$ws.Range("J133").Value = 44563.355
$ws.Range("L133").Value = 44563.355
$ws.Range("N133").Value = -49623.355

$ws.Range("H139").Value = 74333.336
$ws.Range("J139").Value = 74333.336
$ws.Range("L139").Value = 74333.336
$ws.Range("N139").Value = -84613.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 46390
$ws.Range("J133").Value = 46390
$ws.Range("L133").Value = 46390
$ws.Range("N133").Value = -56510

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1381
$ws.Range("I31").Value = 931.7273
$ws.Range("J31").Value = 2113.1482
$ws.Range("K31").Value = 931.7273
$ws.Range("L31").Value = 2113.1482
$ws.Range("M31").Value = -636.7273
$ws.Range("N31").Value = -2703.1482

$ws.Range("H34").Value = 1381
$ws.Range("I34").Value = 931.7273
$ws.Range("J34").Value = 2113.1482
$ws.Range("K34").Value = 931.7273
$ws.Range("L34").Value = 2113.1482
$ws.Range("M34").Value = -729.7273
$ws.Range("N34").Value = -2517.1482

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 96674.336
$ws.Range("I3").Value = 94995
$ws.Range("J3").Value = 100033
$ws.Range("K3").Value = 284985
$ws.Range("L3").Value = 300099
$ws.Range("M3").Value = -284873
$ws.Range("N3").Value = -300323

$ws.Range("H26").Value = 275
$ws.Range("I26").Value = 50
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 150
$ws.Range("L26").Value = 1500
$ws.Range("M26").Value = 138
$ws.Range("N26").Value = -2076

$ws.Range("H32").Value = 421.66666
$ws.Range("J32").Value = 446
$ws.Range("L32").Value = 1338
$ws.Range("N32").Value = -1904

$ws.Range("H33").Value = 281.0909
$ws.Range("I33").Value = 115.71429
$ws.Range("J33").Value = 570.5
$ws.Range("K33").Value = 694.28574
$ws.Range("L33").Value = 3423
$ws.Range("M33").Value = -411.28574
$ws.Range("N33").Value = -3989

$ws.Range("H38").Value = 132.14815
$ws.Range("I38").Value = 180.81818
$ws.Range("J38").Value = 98.6875
$ws.Range("K38").Value = 542.4545400000001
$ws.Range("L38").Value = 296.0625
$ws.Range("M38").Value = -195.4545400000001
$ws.Range("N38").Value = -990.0625

$ws.Range("H41").Value = 777.7778
$ws.Range("I41").Value = 500
$ws.Range("J41").Value = 812.5
$ws.Range("K41").Value = 1500
$ws.Range("L41").Value = 2437.5
$ws.Range("M41").Value = -1162
$ws.Range("N41").Value = -3113.5

$ws.Range("H136").Value = 6312.0454
$ws.Range("I136").Value = 2005
$ws.Range("J136").Value = 7927.1875
$ws.Range("K136").Value = 6015
$ws.Range("L136").Value = 23781.5625
$ws.Range("M136").Value = -915
$ws.Range("N136").Value = -33981.5625

$ws.Range("H139").Value = 1651.8182
$ws.Range("I139").Value = 1651.8182
$ws.Range("K139").Value = 4955.4546
$ws.Range("M139").Value = 184.5454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 627.2820400000001
$ws.Range("I97").Value = 535.4091
$ws.Range("J97").Value = 746.17645
$ws.Range("K97").Value = 535.4091
$ws.Range("L97").Value = 746.17645
$ws.Range("M97").Value = -39.40909999999997
$ws.Range("N97").Value = -1738.17645

$ws.Range("H113").Value = 1533.8334
$ws.Range("I113").Value = 1597
$ws.Range("J113").Value = 1515.7858
$ws.Range("K113").Value = 1597
$ws.Range("L113").Value = 1515.7858
$ws.Range("M113").Value = 573
$ws.Range("N113").Value = -5855.7858

$ws.Range("H138").Value = 63500
$ws.Range("J138").Value = 63500
$ws.Range("L138").Value = 63500
$ws.Range("N138").Value = -73780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3796.7026
$ws.Range("I132").Value = 2843.04
$ws.Range("J132").Value = 5783.5
$ws.Range("K132").Value = 8529.119999999999
$ws.Range("L132").Value = 17350.5
$ws.Range("M132").Value = -5999.119999999999
$ws.Range("N132").Value = -22410.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 44465890
$ws.Range("I2").Value = 57150572
$ws.Range("J2").Value = 69496.5
$ws.Range("K2").Value = 57150572
$ws.Range("L2").Value = 69496.5
$ws.Range("M2").Value = -57150460
$ws.Range("N2").Value = -69720.5

$ws.Range("H24").Value = 16677872
$ws.Range("J24").Value = 13446
$ws.Range("L24").Value = 13446
$ws.Range("N24").Value = -13906
